$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Apply cell formatting (fill/border/font/alignment) by copying
# formats from existing styled cells, so existing style indices are reused
# wherever possible instead of duplicating style entries. ---

$ws.Range("A1").Copy() | Out-Null
foreach ($addr in @("E14","F14","G14","H14","I14","J14","K14","L14","M14","N14","O14","P14","Q14","R14","S14")) { $ws.Range($addr).PasteSpecial(-4122) | Out-Null }

$ws.Range("A2").Copy() | Out-Null
foreach ($addr in @("D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24")) { $ws.Range($addr).PasteSpecial(-4122) | Out-Null }

$ws.Range("B2").Copy() | Out-Null
foreach ($addr in @("E15","F15","G15","H15","K15","L15","O15","P15","Q15","S15","E16","F16","G16","H16","K16","L16","M16","N16","O16","Q16","S16","E17","F17","G17","H17","M17","N17","O17","P17","Q17","S17","E18","F18","G18","H18","K18","L18","M18","N18","P18","Q18","S18","E19","F19","G19","H19","K19","L19","M19","N19","O19","P19","Q19","S19","E20","F20","G20","H20","K20","L20","O20","P20","Q20","S20","E21","F21","G21","H21","N21","O21","Q21","S21","E22","F22","G22","H22","Q22","S22","E23","F23","G23","H23","Q23","S23","E24","F24","G24","H24","K24","L24","M24","N24","P24","Q24","S24")) { $ws.Range($addr).PasteSpecial(-4122) | Out-Null }

$ws.Range("F2").Copy() | Out-Null
foreach ($addr in @("I15","J15","R15","I16","J16","R16","I17","J17","K17","L17","R17","R18","I19","J19","R19","R20","I21","J21","R21","I22","J22","K22","L22","R22","I23","J23","K23","L23","R23","R24")) { $ws.Range($addr).PasteSpecial(-4122) | Out-Null }

$ws.Range("J2").Copy() | Out-Null
foreach ($addr in @("M15","N15","P16","I18","J18","O18","M20","N20","K21","L21","M21","P21","M22","N22","O22","P22","M23","N23","O23","P23","I24","J24","O24")) { $ws.Range($addr).PasteSpecial(-4122) | Out-Null }

# New fill color FFFFF2CC (light yellow) - style 7 (copy base look from F4, then recolor)
$ws.Range("F4").Copy() | Out-Null
foreach ($addr in @("I20","J20")) { $ws.Range($addr).PasteSpecial(-4122) | Out-Null }

$ws.Range("I20:J20").Interior.Color = 13431551

# --- Step 2: Set values ---

$ws.Range("D14").Value = "Altona 5s"

$ws.Range("E14").Value = "2 - Point FGA"
$ws.Range("F14").Value = "2 - Point FGM"
$ws.Range("G14").Value = "3 - Point FGA"
$ws.Range("H14").Value = "3 - Point FGM"
$ws.Range("I14").Value = "Total FGA"
$ws.Range("J14").Value = "Total FGM"
$ws.Range("K14").Value = "FTA"
$ws.Range("L14").Value = "FTM"
$ws.Range("M14").Value = "O - Boards"
$ws.Range("N14").Value = "D - Boards"
$ws.Range("O14").Value = "Assists"
$ws.Range("P14").Value = "Steals"
$ws.Range("Q14").Value = "Blocks"
$ws.Range("R14").Value = "Turnovers"
$ws.Range("S14").Value = "Fouls"

$ws.Range("D15").Value = "Abeal"
$ws.Range("D16").Value = "Ayden"
$ws.Range("D17").Value = "Dayne"
$ws.Range("D18").Value = "Evan"
$ws.Range("D19").Value = "Louis.S"
$ws.Range("D20").Value = "Louis.VB"
$ws.Range("D21").Value = "Lucas"
$ws.Range("D22").Value = "Max"
$ws.Range("D23").Value = "Rikin"
$ws.Range("D24").Value = "Tristyn"

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 3
$ws.Range("N15").Value = 2
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 1
$ws.Range("S15").Value = 0
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 2
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 0
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 1
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 3
$ws.Range("S17").Value = 3
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 4
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 2
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 3
$ws.Range("S18").Value = 0
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 4
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 1
$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = 2
$ws.Range("S19").Value = 2
$ws.Range("E20").Value = 6
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 3
$ws.Range("N20").Value = 2
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 1
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = 1
$ws.Range("S20").Value = 3
$ws.Range("E21").Value = 6
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 6
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 2
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 2
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = 4
$ws.Range("S21").Value = 2
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 9
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 6
$ws.Range("L22").Value = 2
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 2
$ws.Range("O22").Value = 2
$ws.Range("P22").Value = 3
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 4
$ws.Range("E23").Value = 6
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 7
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 5
$ws.Range("L23").Value = 2
$ws.Range("M23").Value = 2
$ws.Range("N23").Value = 2
$ws.Range("O23").Value = 5
$ws.Range("P23").Value = 2
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = 1
$ws.Range("S23").Value = 2
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 4
$ws.Range("J24").Value = 3
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 1
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 2
$ws.Range("P24").Value = 1
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = 1
$ws.Range("S24").Value = 0

# --- Step 3: Row heights (approximate thick-bottom-border row sizing) ---
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 26.5
$ws.Rows.Item(15).RowHeight = 15
$ws.Rows.Item(16).RowHeight = 15
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 15
$ws.Rows.Item(19).RowHeight = 15
$ws.Rows.Item(20).RowHeight = 15
$ws.Rows.Item(21).RowHeight = 15
$ws.Rows.Item(22).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 15
$ws.Rows.Item(24).RowHeight = 15

# --- Step 4: Selection ---
$ws.Range("F11").Select() | Out-Null

